# cv121012a.xlsx — "correção nos dados e inicio da analise PNAD 2009"
#
# The sheet was exported from a pandas multi-index table and carried two
# stray header/section rows that had no data of their own:
#   row 5 -> "situação do domicílio"   (section header for urbana/rural)
#   row 8 (after the first row is gone, this is the original row 8)
#         -> "grandes regiões e unidades da federação" (section header for the UFs)
# Both are removed (EntireRow delete, shifting everything below them up),
# which is why every region row after them moves up by one, and the sheet
# ends up 2 rows shorter (H40 -> H38). The B2 header label for the first
# numeric column is also corrected from the leftover pandas placeholder
# "unnamed: 1_level_1" to "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the leftover pandas column header text.
$ws.Range("B2").Value2 = "total"

# Remove the "situação do domicílio" section-header row (row 5) ...
$ws.Rows(5).Delete()

# ... then remove the "grandes regiões e unidades da federação" section-header
# row, which after the first delete has shifted up into row 7.
$ws.Rows(7).Delete()
